$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid and Absent counts set to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count and Real counts set to 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Rows 5-18: Absent count set to 1
for ($r = 5; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
